$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the anchor paragraph: "Để lấy thay đổi: git pull origin <tên nhánh>"
# (the paragraph that ends with "origin <tên nhánh>")
# ------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("origin <tên nhánh>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $findRng.Paragraphs(1)
$anchorStart = $anchorPara.Range.Start
$anchorEnd = $anchorPara.Range.End
$anchorLen = $anchorEnd - $anchorStart

# ------------------------------------------------------------------
# Duplicate that whole paragraph (text + formatting + its own paragraph
# mark) twice right after itself. Using Copy/Paste (instead of typing
# fresh text) means the new paragraphs inherit the exact same run
# formatting pattern (non-bold lead-in + bold command), including the
# absence of <w:bCs/> on the non-bold run, which matching the target.
# ------------------------------------------------------------------
$anchorPara.Range.Copy()

$p1Start = $anchorEnd
$d.Range($p1Start, $p1Start).Paste()
$p1End = $p1Start + $anchorLen

$p2Start = $p1End
$d.Range($p2Start, $p2Start).Paste()
$p2End = $p2Start + $anchorLen

# ------------------------------------------------------------------
# New paragraph 1: "Để xóa file: git rm <tên file> (sau đó nhớ commit comment)"
# ------------------------------------------------------------------
$oldPrefix = "Để lấy thay đổi: "
$oldCommand = "git pull origin <tên nhánh>"

$newPrefix1 = "Để xóa file: "
$newCommand1 = "git rm <tên file> (sau đó nhớ commit comment)"

$rngA = $d.Range($p1Start, $p1End)
$rngA.Find.Execute($oldPrefix, $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rngA.Text = $newPrefix1
$p1End = $p1End + ($newPrefix1.Length - $oldPrefix.Length)

$rngB = $d.Range($p1Start, $p1End)
$rngB.Find.Execute($oldCommand, $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rngB.Text = $newCommand1
$p1End = $p1End + ($newCommand1.Length - $oldCommand.Length)

# ------------------------------------------------------------------
# New paragraph 2: "Để đổi tên file: git mv <tên_cũ> <tên_mới> (sau đó nhớ commit comment)"
# ------------------------------------------------------------------
$newPrefix2 = "Để đổi tên file:"
$newCommand2 = " git mv <tên_cũ> <tên_mới> (sau đó nhớ commit comment)"

$rngC = $d.Range($p2Start, $p2End)
$rngC.Find.Execute($oldPrefix, $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rngC.Text = $newPrefix2
$p2End = $p2End + ($newPrefix2.Length - $oldPrefix.Length)

$rngD = $d.Range($p2Start, $p2End)
$rngD.Find.Execute($oldCommand, $true, $false, $false, $false, $false, $true, 0, $false, "", 0) | Out-Null
$rngD.Text = $newCommand2
$p2End = $p2End + ($newCommand2.Length - $oldCommand.Length)
